$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the duration (minutes) value in E6; dependent formula cells
# (F6, J6, K6, L6, M6) recalc automatically.
$ws.Range("E6").Value = 2

$excel.Calculate()
